# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on the active worksheet to match the latest scrape.
# Values that look like plain numbers are written with a leading
# apostrophe to force text storage (matching the source data, which
# uses '.' as a thousands separator), then the style is reset to
# Normal so no extra text-format style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "41.671.10"
$ws.Cells.Item(2, 5).Value = "  +0.34%  "
$ws.Cells.Item(3, 4).Value = "2.470.36"
$ws.Cells.Item(3, 5).Value = "  -0.98%  "
$ws.Cells.Item(4, 5).Value = "  +0.24%  "
$ws.Cells.Item(5, 4).Value = "'316.92"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.03%  "
$ws.Cells.Item(6, 4).Value = "'92.69"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.91%  "
$ws.Cells.Item(7, 5).Value = "  +0.77%  "
$ws.Cells.Item(8, 5).Value = "  +0.11%  "
$ws.Cells.Item(9, 4).Value = "'0.514"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.59%  "
$ws.Cells.Item(10, 4).Value = "'0.0896"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +13.93%  "
$ws.Cells.Item(11, 4).Value = "'32.82"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.28%  "
$ws.Cells.Item(12, 5).Value = "  -0.26%  "
$ws.Cells.Item(13, 4).Value = "2.853.12"
$ws.Cells.Item(13, 5).Value = "  -0.87%  "
$ws.Cells.Item(14, 4).Value = "'6.90"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.12%  "
$ws.Cells.Item(15, 4).Value = "'15.72"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.67%  "
$ws.Cells.Item(16, 4).Value = "2.470.28"
$ws.Cells.Item(16, 5).Value = "  -1.19%  "
$ws.Cells.Item(17, 4).Value = "'0.784"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +3.18%  "
$ws.Cells.Item(18, 4).Value = "41.670.17"
$ws.Cells.Item(18, 5).Value = "  +0.22%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0978"
$ws.Cells.Item(19, 5).Value = "  +5.15%  "
$ws.Cells.Item(20, 5).Value = "  +2.16%  "
$ws.Cells.Item(21, 4).Value = "'71.26"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.17%  "
$ws.Cells.Item(22, 4).Value = "'11.44"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.80%  "
$ws.Cells.Item(23, 4).Value = "'239.56"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.33%  "
$ws.Cells.Item(24, 4).Value = "'2.72"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.10%  "
$ws.Cells.Item(25, 4).Value = "'1.91"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.62%  "
$ws.Cells.Item(26, 5).Value = "  -0.05%  "
$ws.Cells.Item(27, 4).Value = "'24.71"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.60%  "
$ws.Cells.Item(28, 5).Value = "  +1.80%  "
$ws.Cells.Item(29, 5).Value = "  +1.24%  "
$ws.Cells.Item(30, 4).Value = "'35.26"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.71%  "
$ws.Cells.Item(31, 4).Value = "'156.03"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.26%  "
$ws.Cells.Item(32, 5).Value = "  +0.60%  "
$ws.Cells.Item(33, 4).Value = "'0.0768"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.14%  "
$ws.Cells.Item(34, 5).Value = "  +0.19%  "
$ws.Cells.Item(35, 4).Value = "'2.51"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.40%  "
$ws.Cells.Item(36, 4).Value = "'17.48"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -2.94%  "
$ws.Cells.Item(37, 4).Value = "'2.88"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.85%  "
$ws.Cells.Item(38, 5).Value = "  +1.13%  "
$ws.Cells.Item(39, 5).Value = "  -2.84%  "
$ws.Cells.Item(40, 5).Value = "  -1.89%  "
$ws.Cells.Item(41, 4).Value = "'3.99"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -3.63%  "
$ws.Cells.Item(42, 5).Value = "  +0.16%  "
$ws.Cells.Item(43, 4).Value = "1.968.16"
$ws.Cells.Item(43, 5).Value = "  -0.26%  "
$ws.Cells.Item(44, 5).Value = "  -0.22%  "
$ws.Cells.Item(45, 4).Value = "'18.85"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -6.45%  "
$ws.Cells.Item(46, 5).Value = "  -1.71%  "
$ws.Cells.Item(47, 4).Value = "'9.08"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.02%  "
$ws.Cells.Item(48, 4).Value = "2.707.42"
$ws.Cells.Item(48, 5).Value = "  -0.71%  "
$ws.Cells.Item(49, 4).Value = "'97.34"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.55%  "
$ws.Cells.Item(50, 4).Value = "'66.89"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.85%  "
$ws.Cells.Item(51, 4).Value = "'52.78"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +4.10%  "
